$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was collected; insert it as a new row 6,
# pushing the existing rows 6-30 down to 7-31 (matches the diff: every
# prior row 6..30 now appears, unchanged, one row lower).
$ws.Rows.Item(6).Insert()

$ws.Range("A6").Value = 2
$ws.Range("B6").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C6").Value = "Coquimbo"
$ws.Range("D6").Value = 44490
$ws.Range("E6").Value = 4
$ws.Range("F6").Value = 100112032
$ws.Range("G6").Value = "Zapallo italiano"
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 600
$ws.Range("K6").Value = 13000
$ws.Range("L6").Value = 15000
$ws.Range("M6").Value = 14000
$ws.Range("N6").Value = "$/caja 60 unidades"
$ws.Range("O6").Value = "Provincia de Limarí"
$ws.Range("P6").Value = 233
$ws.Range("Q6").Value = 60
$ws.Range("R6").Value = "Hortaliza"
